# Add three new sprint-task rows (153-155) to the task tracker sheet,
# mirroring the existing S16/G05 "Analytics & Docs" rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        Row = 153
        A = 'S16'
        B = 'G05'
        C = 'Analytics & Docs'
        D = 'S16_G05_TB002'
        E = 'Implement unified fundamentals ingestion per 16_Unified_Ingestion_Design_v2.md: FundamentalsSnapshotRun lineage table, SymbolResolverService, FundamentalsIngestionService.ingest_screener_csv(), and a manual script that reads backend/data/fundamentals/*.csv and upserts FundamentalsSnapshot.'
        G = 'implemented'
        H = 'Implemented FundamentalsIngestionService with Screener CSV upsert, FundamentalsSnapshotRun lineage table, and CLI script backend/scripts/ingest_screener_fundamentals.py reading backend/data/fundamentals/*.csv.'
    },
    @{
        Row = 154
        A = 'S16'
        B = 'G05'
        C = 'Analytics & Docs'
        D = 'S16_G05_TB003'
        E = 'Implement FactorRiskRebuildService and /api/v1/factors/rebuild endpoint (and optional CLI wrapper) to recompute factor_exposures, risk_model and covariance_matrices for a universe/as_of_date, including basic price-coverage diagnostics.'
        G = 'implemented'
        H = 'Implemented FactorRiskRebuildService and /api/v1/factors/rebuild endpoint to recompute factor_exposures, risk_model, and covariance_matrices for a universe/as_of_date, with basic price coverage diagnostics.'
    },
    @{
        Row = 155
        A = 'S16'
        B = 'G05'
        C = 'Analytics & Docs'
        D = 'S16_G05_TD002'
        E = 'Finalize and maintain unified ingestion docs (16_Unified_Ingestion_Design.md + 16_Unified_Ingestion_Design_v2.md) so they stay aligned with the implemented fundamentals ingestion and factor/risk rebuild flows.'
        G = 'implemented'
        H = 'Finalised unified ingestion documentation (Files 16 and 16_v2) to match the implemented fundamentals ingestion and factor/risk rebuild flows.'
    }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.A
    $ws.Cells.Item($r.Row, 2).Value2 = $r.B
    $ws.Cells.Item($r.Row, 3).Value2 = $r.C
    $ws.Cells.Item($r.Row, 4).Value2 = $r.D
    $ws.Cells.Item($r.Row, 5).Value2 = $r.E
    # Column F ("remarks") is intentionally left blank, matching rows 150-152.
    $ws.Cells.Item($r.Row, 7).Value2 = $r.G
    $ws.Cells.Item($r.Row, 8).Value2 = $r.H
}
